$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]1.02
$ws.Range("C2").Value = [double]1.024422830022724
$ws.Range("D2").Value = [double]1.027588801413152
$ws.Range("E2").Value = [double]0.9926147277508489
$ws.Range("F2").Value = [double]1.022904906233393
$ws.Range("I2").Value = [double]1.028489610443474
$ws.Range("J2").Value = [double]1.029597529870352
$ws.Range("K2").Value = [double]1.030408090989892
$ws.Range("L2").Value = [double]0.9955398523336033
$ws.Range("M2").Value = [double]1.025737912806098
$ws.Range("N2").Value = [double]1.031059676459259
$ws.Range("B3").Value = [double]1.02
$ws.Range("C3").Value = [double]1.026212814749323
$ws.Range("D3").Value = [double]1.028872378576037
$ws.Range("E3").Value = [double]0.9936372048519304
$ws.Range("F3").Value = [double]1.025344888860022
$ws.Range("I3").Value = [double]1.028843333467778
$ws.Range("J3").Value = [double]1.031022072208288
$ws.Range("K3").Value = [double]1.031498139546625
$ws.Range("L3").Value = [double]0.9963617723202692
$ws.Range("M3").Value = [double]1.027980210307358
$ws.Range("N3").Value = [double]1.032486241810709
$ws.Range("B4").Value = [double]1.02
$ws.Range("C4").Value = [double]1.027366031667895
$ws.Range("D4").Value = [double]1.029698729517703
$ws.Range("E4").Value = [double]0.9942998659930995
$ws.Range("F4").Value = [double]1.026917845032661
$ws.Range("I4").Value = [double]1.029069123956171
$ws.Range("J4").Value = [double]1.031938641591444
$ws.Range("K4").Value = [double]1.032198734627715
$ws.Range("L4").Value = [double]0.9968940712668345
$ws.Range("M4").Value = [double]1.029424996647089
$ws.Range("N4").Value = [double]1.03340411282752
$ws.Range("B5").Value = [double]1.02
$ws.Range("C5").Value = [double]1.027849665623369
$ws.Range("D5").Value = [double]1.030045135706145
$ws.Range("E5").Value = [double]0.9945786998346017
$ws.Range("F5").Value = [double]1.027577747262248
$ws.Range("I5").Value = [double]1.029163311877925
$ws.Range("J5").Value = [double]1.032322741345528
$ws.Range("K5").Value = [double]1.032492145381884
$ws.Range("L5").Value = [double]0.997117960005301
$ws.Range("M5").Value = [double]1.030030952622353
$ws.Range("N5").Value = [double]1.033788758047312
$ws.Range("B6").Value = [double]1.02
$ws.Range("C6").Value = [double]1.027930801481048
$ws.Range("D6").Value = [double]1.030103241132478
$ws.Range("E6").Value = [double]0.9946255319796338
$ws.Range("F6").Value = [double]1.027688468595669
$ws.Range("I6").Value = [double]1.029179083546305
$ws.Range("J6").Value = [double]1.032387162066416
$ws.Range("K6").Value = [double]1.032541345131639
$ws.Range("L6").Value = [double]0.9971555583673453
$ws.Range("M6").Value = [double]1.030132612329695
$ws.Range("N6").Value = [double]1.033853270253013
$ws.Range("B7").Value = [double]1.02
$ws.Range("C7").Value = [double]1.027372498606491
$ws.Range("D7").Value = [double]1.029703362089424
$ws.Range("E7").Value = [double]0.9943035907978915
$ws.Range("F7").Value = [double]1.026926668004861
$ws.Range("I7").Value = [double]1.029070385377577
$ws.Range("J7").Value = [double]1.031943778741545
$ws.Range("K7").Value = [double]1.03220265957462
$ws.Range("L7").Value = [double]0.9968970624459041
$ws.Range("M7").Value = [double]1.029433099040615
$ws.Range("N7").Value = [double]1.033409257272963
$ws.Range("B8").Value = [double]1.02
$ws.Range("C8").Value = [double]1.025028818339592
$ws.Range("D8").Value = [double]1.028023473572437
$ws.Range("E8").Value = [double]0.9929600610674294
$ws.Range("F8").Value = [double]1.023730749495833
$ws.Range("I8").Value = [double]1.028609796537922
$ws.Range("J8").Value = [double]1.030080050436155
$ws.Range("K8").Value = [double]1.030777468482558
$ws.Range("L8").Value = [double]0.9958175282591053
$ws.Range("M8").Value = [double]1.026496998441038
$ws.Range("N8").Value = [double]1.031542882259608
$ws.Range("B9").Value = [double]1.02
$ws.Range("C9").Value = [double]1.020859394789898
$ws.Range("D9").Value = [double]1.025030332405298
$ws.Range("E9").Value = [double]0.9906006454969559
$ws.Range("F9").Value = [double]1.018052404185311
$ws.Range("I9").Value = [double]1.027774234855501
$ws.Range("J9").Value = [double]1.02675515300316
$ws.Range("K9").Value = [double]1.028229114859574
$ws.Range("L9").Value = [double]0.9939188001724441
$ws.Range("M9").Value = [double]1.021274662798654
$ws.Range("N9").Value = [double]1.028213263090887
$ws.Range("B10").Value = [double]1.02
$ws.Range("C10").Value = [double]1.018051643037462
$ws.Range("D10").Value = [double]1.023011718088198
$ws.Range("E10").Value = [double]0.989033133672735
$ws.Range("F10").Value = [double]1.014232986650913
$ws.Range("I10").Value = [double]1.027200719498039
$ws.Range("J10").Value = [double]1.02450986115748
$ws.Range("K10").Value = [double]1.0265043732851
$ws.Range("L10").Value = [double]0.9926553831429383
$ws.Range("M10").Value = [double]1.017758202523089
$ws.Range("N10").Value = [double]1.025964782673247
$ws.Range("B11").Value = [double]1.02
$ws.Range("C11").Value = [double]1.016828835561598
$ws.Range("D11").Value = [double]1.022131904967933
$ws.Range("E11").Value = [double]0.988355674866747
$ws.Range("F11").Value = [double]1.012570558740999
$ws.Range("I11").Value = [double]1.026948388725683
$ws.Range("J11").Value = [double]1.023530530249539
$ws.Range("K11").Value = [double]1.025751192905727
$ws.Range("L11").Value = [double]0.9921088820399291
$ws.Range("M11").Value = [double]1.016226749161528
$ws.Range("N11").Value = [double]1.024984061003086
$ws.Range("B12").Value = [double]1.02
$ws.Range("C12").Value = [double]1.016373543750912
$ws.Range("D12").Value = [double]1.021804220705268
$ws.Range("E12").Value = [double]0.9881042295826724
$ws.Range("F12").Value = [double]1.011951720118581
$ws.Range("I12").Value = [double]1.02685405423187
$ws.Range("J12").Value = [double]1.023165670579036
$ws.Range("K12").Value = [double]1.025470454132037
$ws.Range("L12").Value = [double]0.9919059725120875
$ws.Range("M12").Value = [double]1.0156565324943
$ws.Range("N12").Value = [double]1.02461868319
$ws.Range("B13").Value = [double]1.02
$ws.Range("C13").Value = [double]1.016471255052983
$ws.Range("D13").Value = [double]1.02187455035416
$ws.Range("E13").Value = [double]0.9881581567098651
$ws.Range("F13").Value = [double]1.0120845244637
$ws.Range("I13").Value = [double]1.026874316917392
$ws.Range("J13").Value = [double]1.023243984066479
$ws.Range("K13").Value = [double]1.025530717940544
$ws.Range("L13").Value = [double]0.9919494934313052
$ws.Range("M13").Value = [double]1.015778908468545
$ws.Range("N13").Value = [double]1.024697107891579
$ws.Range("B14").Value = [double]1.02
$ws.Range("C14").Value = [double]1.016791223342367
$ws.Range("D14").Value = [double]1.022104836605126
$ws.Range("E14").Value = [double]0.9883348863814464
$ws.Range("F14").Value = [double]1.012519432938491
$ws.Range("I14").Value = [double]1.026940603445543
$ws.Range("J14").Value = [double]1.023500393270563
$ws.Range("K14").Value = [double]1.025728006933345
$ws.Range("L14").Value = [double]0.9920921077337197
$ws.Range("M14").Value = [double]1.016179643012606
$ws.Range("N14").Value = [double]1.024953881226144
$ws.Range("B15").Value = [double]1.02
$ws.Range("C15").Value = [double]1.016988221461393
$ws.Range("D15").Value = [double]1.022246606016216
$ws.Range("E15").Value = [double]0.9884438009545853
$ws.Range("F15").Value = [double]1.012787215551051
$ws.Range("I15").Value = [double]1.026981364024371
$ws.Range("J15").Value = [double]1.023658229907642
$ws.Range("K15").Value = [double]1.025849433544034
$ws.Range("L15").Value = [double]0.9921799884222134
$ws.Range("M15").Value = [double]1.016426366354947
$ws.Range("N15").Value = [double]1.02511194200935
$ws.Range("B16").Value = [double]1.02
$ws.Range("C16").Value = [double]1.018132646017581
$ws.Range("D16").Value = [double]1.02306998574418
$ws.Range("E16").Value = [double]0.9890781214508737
$ws.Range("F16").Value = [double]1.014343131388584
$ws.Range("I16").Value = [double]1.02721738106261
$ws.Range("J16").Value = [double]1.02457470429851
$ws.Range("K16").Value = [double]1.026554223848686
$ws.Range("L16").Value = [double]0.9926916645766087
$ws.Range("M16").Value = [double]1.017859650820513
$ws.Range("N16").Value = [double]1.026029717898975
$ws.Range("B17").Value = [double]1.02
$ws.Range("C17").Value = [double]1.018848609419334
$ws.Range("D17").Value = [double]1.023584919270207
$ws.Range("E17").Value = [double]0.9894763578477731
$ws.Range("F17").Value = [double]1.015316781854179
$ws.Range("I17").Value = [double]1.027364353581035
$ws.Range("J17").Value = [double]1.025147664710534
$ws.Range("K17").Value = [double]1.026994605125327
$ws.Range("L17").Value = [double]0.9930127773692701
$ws.Range("M17").Value = [double]1.018756324587293
$ws.Range("N17").Value = [double]1.026603491980504
$ws.Range("B18").Value = [double]1.02
$ws.Range("C18").Value = [double]1.01926554229084
$ws.Range("D18").Value = [double]1.023884718856343
$ws.Range("E18").Value = [double]0.9897087662937551
$ws.Range("F18").Value = [double]1.015883869669331
$ws.Range("I18").Value = [double]1.027449695229591
$ws.Range("J18").Value = [double]1.025481179044447
$ws.Range("K18").Value = [double]1.027250860132094
$ws.Range("L18").Value = [double]0.9932001317071766
$ws.Range("M18").Value = [double]1.019278492899053
$ws.Range("N18").Value = [double]1.026937479943026
$ws.Range("B19").Value = [double]1.02
$ws.Range("C19").Value = [double]1.019407591788718
$ws.Range("D19").Value = [double]1.023986849670841
$ws.Range("E19").Value = [double]0.9897880325774039
$ws.Range("F19").Value = [double]1.016077093278413
$ws.Range("I19").Value = [double]1.027478729457413
$ws.Range("J19").Value = [double]1.02559478358724
$ws.Range("K19").Value = [double]1.027338133245428
$ws.Range("L19").Value = [double]0.993264023964098
$ws.Range("M19").Value = [double]1.01945639657296
$ws.Range("N19").Value = [double]1.0270512458173
$ws.Range("B20").Value = [double]1.02
$ws.Range("C20").Value = [double]1.018771863508937
$ws.Range("D20").Value = [double]1.023529729061246
$ws.Range("E20").Value = [double]0.9894336180355766
$ws.Range("F20").Value = [double]1.015212404155686
$ws.Range("I20").Value = [double]1.02734862468978
$ws.Range("J20").Value = [double]1.025086262388301
$ws.Range("K20").Value = [double]1.02694741981565
$ws.Range("L20").Value = [double]0.9929783193490043
$ws.Range("M20").Value = [double]1.018660207834988
$ws.Range("N20").Value = [double]1.026542002459931
$ws.Range("B21").Value = [double]1.02
$ws.Range("C21").Value = [double]1.016697030898014
$ws.Range("D21").Value = [double]1.022037047568816
$ws.Range("E21").Value = [double]0.9882828385668255
$ws.Range("F21").Value = [double]1.012391400578396
$ws.Range("I21").Value = [double]1.026921100548515
$ws.Range("J21").Value = [double]1.023424917515246
$ws.Range("K21").Value = [double]1.025669937286123
$ws.Range("L21").Value = [double]0.9920501090198107
$ws.Range("M21").Value = [double]1.016061674763669
$ws.Range("N21").Value = [double]1.024878298286597
$ws.Range("B22").Value = [double]1.02
$ws.Range("C22").Value = [double]1.015386198853587
$ws.Range("D22").Value = [double]1.021093423371851
$ws.Range("E22").Value = [double]0.9875604150241496
$ws.Range("F22").Value = [double]1.010609953141295
$ws.Range("I22").Value = [double]1.026648779927287
$ws.Range("J22").Value = [double]1.022374030012878
$ws.Range("K22").Value = [double]1.024861088056267
$ws.Range("L22").Value = [double]0.991467000034148
$ws.Range("M22").Value = [double]1.014419944121872
$ws.Range("N22").Value = [double]1.023825918403436
$ws.Range("B23").Value = [double]1.02
$ws.Range("C23").Value = [double]1.01608170331368
$ws.Range("D23").Value = [double]1.021594148183446
$ws.Range("E23").Value = [double]0.9879432794636459
$ws.Range("F23").Value = [double]1.011555084800407
$ws.Range("I23").Value = [double]1.026793478339864
$ws.Range("J23").Value = [double]1.022931734272548
$ws.Range("K23").Value = [double]1.025290415923805
$ws.Range("L23").Value = [double]0.9917760702887607
$ws.Range("M23").Value = [double]1.015291023201371
$ws.Range("N23").Value = [double]1.024384414667125
$ws.Range("B24").Value = [double]1.02
$ws.Range("C24").Value = [double]1.018806543771555
$ws.Range("D24").Value = [double]1.023554668847396
$ws.Range("E24").Value = [double]0.9894529299347241
$ws.Range("F24").Value = [double]1.015259570495165
$ws.Range("I24").Value = [double]1.027355733087894
$ws.Range("J24").Value = [double]1.025114009566658
$ws.Range("K24").Value = [double]1.026968742716263
$ws.Range("L24").Value = [double]0.9929938892766438
$ws.Range("M24").Value = [double]1.018703641471348
$ws.Range("N24").Value = [double]1.026569789042463
$ws.Range("B25").Value = [double]1.02
$ws.Range("C25").Value = [double]1.021942135555366
$ws.Range("D25").Value = [double]1.025808141802645
$ws.Range("E25").Value = [double]0.9912096547607046
$ws.Range("F25").Value = [double]1.019526178459521
$ws.Range("I25").Value = [double]1.027993123382206
$ws.Range("J25").Value = [double]1.027619680805744
$ws.Range("K25").Value = [double]1.028892407345687
$ws.Range("L25").Value = [double]0.9944092447426411
$ws.Range("M25").Value = [double]1.022630745666791
$ws.Range("N25").Value = [double]1.029079018622113

Write-Output "Updated 264 vm_pu cells for 380 kV case (B=1.02 slack voltage)"